# Auto update Excel log
# Appends newly-logged sensor readings to the Humidity, Temperature,
# Proximity and Camera sheets of the SeniorConnect master log.
#
# The source data stores everything as plain text (dates like
# "2026-02-01" and percentages like "78.4%" are literal strings, not
# Excel date/number values), so the Date column (A) and, where present,
# a percentage-looking Value column (E) are pre-formatted as Text
# before the values are written so Excel does not silently convert
# them into date serials / numeric percentages.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        $Rows,
        [bool]$ValueColumnIsPercent
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    # Keep the Date column as literal text (e.g. "2026-02-01").
    $ws.Range($ws.Cells.Item($StartRow, 1), $ws.Cells.Item($endRow, 1)).NumberFormat = "@"

    # Some sheets store percentage-looking strings (e.g. "78.4%") in
    # column E; keep those as literal text too.
    if ($ValueColumnIsPercent) {
        $ws.Range($ws.Cells.Item($StartRow, 5), $ws.Cells.Item($endRow, 5)).NumberFormat = "@"
    }

    $r = $StartRow
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }
}

# --- Humidity sheet: append rows 179-188 ---
$humidityRows = @(
    ,@("2026-02-01", "14:16:02", "14:00", "Bathroom", "78.4%", "Active")
    ,@("2026-02-01", "14:16:22", "14:00", "Bathroom", "77.1%", "Active")
    ,@("2026-02-01", "14:16:24", "14:00", "Bathroom", "76.9%", "Active")
    ,@("2026-02-01", "14:16:27", "14:00", "Bathroom", "78.0%", "Active")
    ,@("2026-02-01", "14:16:33", "14:00", "Bathroom", "76.9%", "Active")
    ,@("2026-02-01", "14:16:37", "14:00", "Bathroom", "77.8%", "Active")
    ,@("2026-02-01", "14:16:42", "14:00", "Bathroom", "76.9%", "Active")
    ,@("2026-02-01", "14:16:54", "14:00", "Bathroom", "77.8%", "Active")
    ,@("2026-02-01", "14:16:56", "14:00", "Bathroom", "76.9%", "Active")
    ,@("2026-02-01", "14:16:57", "14:00", "Bathroom", "77.9%", "Active")
)
Append-Rows "Humidity" 179 $humidityRows $true

# --- Temperature sheet: append rows 100-109 ---
$temperatureRows = @(
    ,@("2026-02-01", "14:16:03", "14:00", "Bathroom", "29.5C", "Active")
    ,@("2026-02-01", "14:16:23", "14:00", "Bathroom", "29.5C", "Active")
    ,@("2026-02-01", "14:16:25", "14:00", "Bathroom", "29.4C", "Active")
    ,@("2026-02-01", "14:16:28", "14:00", "Bathroom", "29.5C", "Active")
    ,@("2026-02-01", "14:16:33", "14:00", "Bathroom", "29.4C", "Active")
    ,@("2026-02-01", "14:16:38", "14:00", "Bathroom", "29.4C", "Active")
    ,@("2026-02-01", "14:16:43", "14:00", "Bathroom", "29.5C", "Active")
    ,@("2026-02-01", "14:16:55", "14:00", "Bathroom", "29.5C", "Active")
    ,@("2026-02-01", "14:16:56", "14:00", "Bathroom", "29.4C", "Active")
    ,@("2026-02-01", "14:16:58", "14:00", "Bathroom", "29.5C", "Active")
)
Append-Rows "Temperature" 100 $temperatureRows $false

# --- Proximity sheet: append rows 29-32 ---
$proximityRows = @(
    ,@("2026-02-01", "14:16:03", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-02-01", "14:16:21", "14:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
    ,@("2026-02-01", "14:16:44", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-02-01", "14:16:54", "14:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
)
Append-Rows "Proximity" 29 $proximityRows $false

# --- Camera sheet: append rows 18-19 ---
$cameraRows = @(
    ,@("2026-02-01", "14:16:21", "14:00", "Living Room Main Door", "Image Captured", "Active")
    ,@("2026-02-01", "14:16:53", "14:00", "Living Room Main Door", "Image Captured", "Active")
)
Append-Rows "Camera" 18 $cameraRows $false
